$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: the "Sex" breakdown values for the "below poverty level" block.
# Rows 4-16 were "male" -> now " Male" (leading space, capitalised).
$ws.Range("D4:D16").Value = " Male"

# Section header labels in column A.
# Row 31 was "Male:" -> now " Male:" (leading space added).
$ws.Range("A31").Value = " Male:"
# Row 45 was "Female:" -> now "Fe Male:".
$ws.Range("A45").Value = "Fe Male:"

# Rows 17-29 were "female" -> now "Female" (capitalised, no trailing colon).
$ws.Range("D17:D29").Value = "Female"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("D31").Select()
